# Planning.xlsx update:
#  - Analyse/tijdschema bijgewerkt (tijden verschoven)
#  - "ERD + uitloop" opgesplitst in aparte taken "ERD" en "Uitloop"
#  - "Testplan" hernoemd/samengevoegd tot "testplan + gesprek"
#  - "strokendiagram" vervangen door "klassendiagram"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (14:00 - 14:15, was 14:00-14:30): end time shortened, task renamed to "Uitloop"
# (set before row 4's "ERD" so the shared-string table gets "Uitloop" first, matching
# the order newly-added strings appear in the saved workbook)
$ws.Range("C5").Value = 0.59375
$ws.Range("E5").Value = "Uitloop"

# Row 4 (13:45 - 14:00): task renamed from "ERD + uitloop" to "ERD"
$ws.Range("E4").Value = "ERD"

# Row 6 (14:15 - 15:00, was 14:30-15:15): start/end shifted, task renamed to "klassendiagram"
$ws.Range("B6").Value = 0.59375
$ws.Range("C6").Value = 0.625
$ws.Range("E6").Value = "klassendiagram"

# Row 7 (15:00 - 16:00, was 15:15-16:00): start shifted, task renamed to "testplan + gesprek"
$ws.Range("B7").Value = 0.625
$ws.Range("E7").Value = "testplan + gesprek"

# Update selected cell to reflect saved cursor position
$ws.Range("E7").Select()
